$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and G to be stored as text so the numeric-looking
# strings ("283.60", "11", etc.) are preserved verbatim instead of
# being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = '283.60'
$ws.Range("G2").Value = '11'

$ws.Range("D3").Value = '20.98'
$ws.Range("G3").Value = '11'

$ws.Range("D4").Value = '6.220'
$ws.Range("G4").Value = '11'

$ws.Range("D5").Value = '0.06198'
$ws.Range("G5").Value = '11'

$ws.Range("D6").Value = '3.581'
$ws.Range("G6").Value = '11'

$ws.Range("D7").Value = '6.561'
$ws.Range("G7").Value = '11'

$ws.Range("D8").Value = '1.481'
$ws.Range("G8").Value = '11'

$ws.Range("D9").Value = '0.8171'
$ws.Range("G9").Value = '11'

$ws.Range("D10").Value = '0.01389'
$ws.Range("G10").Value = '11'

$ws.Range("D11").Value = '0.1644'
$ws.Range("G11").Value = '11'

$ws.Range("D12").Value = '0.08307'
$ws.Range("G12").Value = '11'

$ws.Range("D13").Value = '0.03611'
$ws.Range("G13").Value = '11'

$ws.Range("D14").Value = '0.03135'
$ws.Range("G14").Value = '11'

$ws.Range("D15").Value = '0.09134'
$ws.Range("G15").Value = '11'

$ws.Range("D16").Value = '3.694'
$ws.Range("G16").Value = '11'

$ws.Range("D17").Value = '0.001637'
$ws.Range("G17").Value = '11'

$ws.Range("D18").Value = '0.04664'
$ws.Range("G18").Value = '11'

$ws.Range("D19").Value = '0.006466'
$ws.Range("G19").Value = '11'

$ws.Range("D20").Value = '0.006197'
$ws.Range("G20").Value = '11'

$ws.Range("D21").Value = '0.001066'
$ws.Range("G21").Value = '11'

$ws.Range("G22").Value = '11'

$ws.Range("D23").Value = '3.819'
$ws.Range("G23").Value = '11'

$ws.Range("D24").Value = '2.323'
$ws.Range("G24").Value = '11'

$ws.Range("D25").Value = '0.3383'
$ws.Range("G25").Value = '11'

$ws.Range("D26").Value = '0.1249'
$ws.Range("G26").Value = '11'

$ws.Range("G27").Value = '11'

$ws.Range("G28").Value = '11'

$ws.Range("G29").Value = '11'

$ws.Range("G30").Value = '11'

$ws.Range("G31").Value = '11'

$ws.Range("G32").Value = '11'

$ws.Range("G33").Value = '11'

$ws.Range("G34").Value = '11'

$ws.Range("G35").Value = '11'

$ws.Range("G36").Value = '11'

$ws.Range("G37").Value = '11'

$ws.Range("G38").Value = '11'

$ws.Range("G39").Value = '11'

$ws.Range("G40").Value = '11'

$ws.Range("D41").Value = '0.007053'
$ws.Range("G41").Value = '11'

$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").Value = '0.1108'
$ws.Range("E42").Value = '41BKEXTokenBKK'
$ws.Range("G42").Value = '11'

$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").Value = '0.003518'
$ws.Range("E43").Value = '42CEJICEJI'
$ws.Range("G43").Value = '11'

$ws.Range("G44").Value = '11'

$ws.Range("D45").Value = '0.00006449'
$ws.Range("G45").Value = '11'

$ws.Range("G46").Value = '11'

$ws.Range("D47").Value = '0.9993'
$ws.Range("G47").Value = '11'

$ws.Range("D48").Value = '0.002733'
$ws.Range("G48").Value = '11'

$ws.Range("D49").Value = '0.00001899'
$ws.Range("G49").Value = '11'

$ws.Range("G50").Value = '11'

$ws.Range("G51").Value = '11'
